$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Factuur")

# Fill row 31 (A31:H31) with letters a-h
$ws.Range("A31").Value = "a"
$ws.Range("B31").Value = "b"
$ws.Range("C31").Value = "c"
$ws.Range("D31").Value = "d"
$ws.Range("E31").Value = "e"
$ws.Range("F31").Value = "f"
$ws.Range("G31").Value = "g"
$ws.Range("H31").Value = "h"

# Update selection to I31 (single cell)
$ws.Range("I31").Select()

# Update workbook view window position
$excel.ActiveWindow.Left = 10275
$excel.ActiveWindow.Top = 1665
